$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.66%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.57%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.203"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.99%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06962"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'4.28%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.419"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.19%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.553"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.01%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.395"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'3.11%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.8992"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-3.96%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1605"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.03%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07540"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'15.81%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07732"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.78%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02937"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.70%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09009"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.38%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001572"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.20%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0006485"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.60%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006537"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.10%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.486"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.14%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.22%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3244"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.32%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1336"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.38%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.052"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.03%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'5.07%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'1.42%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'2.85%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-7.34%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-6.04%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001673"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'3.78%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04366"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.06%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006926"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.43%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1246"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.10%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'2.87%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01179"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.90%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005829"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.86%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E47").Value = "'-0.16%"
$ws.Range("E47").Style = "Normal"
